# Adds season-record columns (Wins, Losses, Ties) to the worksheet.
# Mirrors the behaviour of the new "get season record" downloader code:
# three new trailing columns are appended with the team's Wins/Losses/Ties,
# repeated on every player row, using the same header style as the other
# header cells in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 54

# New header cells (copy style from an existing header cell so the look
# - bold, bordered, centered - matches the rest of row 1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AA1:AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

$wins = 85
$losses = 77
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
